$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# Update the RMA identifiers / related Salesforce Ids for the three
# maintenance-grid rows (new automation run: RMA-MG41-* replacing RMA-8XKB-*).
$ws.Range("E2").Value = "RMA-MG41-001"
$ws.Range("F2").Value = "RMA-MG41-1-1"
$ws.Range("J2").Value = "a7s5f000000xL33AAE"

$ws.Range("E3").Value = "RMA-MG41-002"
$ws.Range("F3").Value = "RMA-MG41-1-2"
$ws.Range("J3").Value = "a7s5f000000xL34AAE"

$ws.Range("E4").Value = "RMA-MG41-003"
$ws.Range("F4").Value = "RMA-MG41-1-3"
$ws.Range("J4").Value = "a7s5f000000xL35AAE"
